# Update "想去人数" (want-to-go headcount) values in column F for a handful
# of rows on the "展览" sheet and the matching rows on the "全部类型" sheet
# (the latter has one extra row inserted above the shifted entries).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 120
$ws1.Range("F3").Value = 221
$ws1.Range("F5").Value = 6624
$ws1.Range("F9").Value = 6030
$ws1.Range("F12").Value = 1239
$ws1.Range("F14").Value = 88
$ws1.Range("F21").Value = 4346
$ws1.Range("F22").Value = 46
$ws1.Range("F23").Value = 15
$ws1.Range("F25").Value = 27

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 120
$ws4.Range("F3").Value = 221
$ws4.Range("F5").Value = 6624
$ws4.Range("F9").Value = 6030
$ws4.Range("F12").Value = 1239
$ws4.Range("F14").Value = 88
$ws4.Range("F21").Value = 4346
$ws4.Range("F23").Value = 46
$ws4.Range("F24").Value = 15
$ws4.Range("F26").Value = 27
